$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 12874987
$ws.Range("I70").Value = 18668082
$ws.Range("J70").Value = 1442.2222
$ws.Range("K70").Value = 56004246
$ws.Range("L70").Value = 4326.6666
$ws.Range("M70").Value = -56003976
$ws.Range("N70").Value = -4866.6666
$ws.Range("H73").Value = 12874987
$ws.Range("I73").Value = 18668082
$ws.Range("J73").Value = 1442.2222
$ws.Range("K73").Value = 56004246
$ws.Range("L73").Value = 4326.6666
$ws.Range("M73").Value = -56003310
$ws.Range("N73").Value = -6198.6666
$ws.Range("H107").Value = 5485.4
$ws.Range("I107").Value = 7074.4
$ws.Range("K107").Value = 7074.4
$ws.Range("M107").Value = -5154.4
$ws.Range("H112").Value = 9083.333000000001
$ws.Range("I112").Value = 600
$ws.Range("J112").Value = 14738.889
$ws.Range("K112").Value = 1800
$ws.Range("L112").Value = 44216.667
$ws.Range("M112").Value = -692
$ws.Range("N112").Value = -46432.667
$ws.Range("H116").Value = 5151.5
$ws.Range("I116").Value = 2962.6365
$ws.Range("J116").Value = 6756.6665
$ws.Range("K116").Value = 2962.6365
$ws.Range("L116").Value = 6756.6665
$ws.Range("M116").Value = 479.3634999999999
$ws.Range("N116").Value = -13640.6665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 603339.6
$ws.Range("I32").Value = 3231.675
$ws.Range("K32").Value = 3231.675
$ws.Range("M32").Value = -2944.675
$ws.Range("H74").Value = 3998.4707
$ws.Range("I74").Value = 866.43475
$ws.Range("J74").Value = 10547.272
$ws.Range("K74").Value = 866.43475
$ws.Range("L74").Value = 10547.272
$ws.Range("M74").Value = 7.565249999999992
$ws.Range("N74").Value = -12295.272
$ws.Range("H77").Value = 3998.4707
$ws.Range("I77").Value = 866.43475
$ws.Range("J77").Value = 10547.272
$ws.Range("K77").Value = 4332.17375
$ws.Range("L77").Value = 52736.36
$ws.Range("M77").Value = 35.82625000000007
$ws.Range("N77").Value = -61472.36

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12755.365
$ws.Range("I31").Value = 1513.4814
$ws.Range("J31").Value = 24896.6
$ws.Range("K31").Value = 1513.4814
$ws.Range("L31").Value = 24896.6
$ws.Range("M31").Value = -1218.4814
$ws.Range("N31").Value = -25486.6
$ws.Range("H34").Value = 12755.365
$ws.Range("I34").Value = 1513.4814
$ws.Range("J34").Value = 24896.6
$ws.Range("K34").Value = 1513.4814
$ws.Range("L34").Value = 24896.6
$ws.Range("M34").Value = -1311.4814
$ws.Range("N34").Value = -25300.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 166.83333
$ws.Range("I13").Value = 173.66667
$ws.Range("J13").Value = 160
$ws.Range("K13").Value = 521.00001
$ws.Range("L13").Value = 480
$ws.Range("M13").Value = -353.00001
$ws.Range("N13").Value = -816
$ws.Range("H47").Value = 288
$ws.Range("I47").Value = 82
$ws.Range("K47").Value = 246
$ws.Range("M47").Value = 185
$ws.Range("H51").Value = 707.1429000000001
$ws.Range("I51").Value = 600
$ws.Range("J51").Value = 750
$ws.Range("K51").Value = 1800
$ws.Range("L51").Value = 2250
$ws.Range("M51").Value = -1340
$ws.Range("N51").Value = -3170
$ws.Range("H55").Value = 682.3077
$ws.Range("I55").Value = 100
$ws.Range("J55").Value = 730.8333
$ws.Range("K55").Value = 300
$ws.Range("L55").Value = 2192.4999
$ws.Range("M55").Value = -123
$ws.Range("N55").Value = -2546.4999
$ws.Range("H57").Value = 2700
$ws.Range("I57").Value = 1500
$ws.Range("K57").Value = 4500
$ws.Range("M57").Value = -3941
$ws.Range("H61").Value = 83.318184
$ws.Range("I61").Value = 53.266666
$ws.Range("K61").Value = 159.799998
$ws.Range("M61").Value = 55.20000199999998
$ws.Range("H63").Value = 6248.9165
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 6248.9165
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 18746.7495
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -20244.7495
$ws.Range("H64").Value = 4100
$ws.Range("I64").Value = 1500
$ws.Range("J64").Value = 4750
$ws.Range("K64").Value = 4500
$ws.Range("L64").Value = 14250
$ws.Range("M64").Value = -4230
$ws.Range("N64").Value = -14790
$ws.Range("H66").Value = 6248.9165
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 6248.9165
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 56240.2485
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -63728.2485
$ws.Range("H67").Value = 4100
$ws.Range("I67").Value = 1500
$ws.Range("J67").Value = 4750
$ws.Range("K67").Value = 4500
$ws.Range("L67").Value = 14250
$ws.Range("M67").Value = -3564
$ws.Range("N67").Value = -16122
$ws.Range("H69").Value = 992.1818
$ws.Range("I69").Value = 300
$ws.Range("J69").Value = 1146
$ws.Range("K69").Value = 900
$ws.Range("L69").Value = 3438
$ws.Range("M69").Value = -89
$ws.Range("N69").Value = -5060
$ws.Range("H72").Value = 992.1818
$ws.Range("I72").Value = 300
$ws.Range("J72").Value = 1146
$ws.Range("K72").Value = 2700
$ws.Range("L72").Value = 10314
$ws.Range("M72").Value = 1356
$ws.Range("N72").Value = -18426
$ws.Range("H76").Value = 2350
$ws.Range("I76").Value = 1900
$ws.Range("K76").Value = 5700
$ws.Range("M76").Value = -5317
$ws.Range("H79").Value = 2350
$ws.Range("I79").Value = 1900
$ws.Range("K79").Value = 5700
$ws.Range("M79").Value = -4374
$ws.Range("H80").Value = 4333.3335
$ws.Range("I80").Value = 1150
$ws.Range("J80").Value = 4731.25
$ws.Range("K80").Value = 3450
$ws.Range("L80").Value = 14193.75
$ws.Range("M80").Value = -2514
$ws.Range("N80").Value = -16065.75
$ws.Range("H81").Value = 2197.6191
$ws.Range("J81").Value = 2271.0527
$ws.Range("L81").Value = 6813.158100000001
$ws.Range("N81").Value = -9059.158100000001
$ws.Range("H82").Value = 4016.6667
$ws.Range("I82").Value = 1900
$ws.Range("J82").Value = 4440
$ws.Range("K82").Value = 5700
$ws.Range("L82").Value = 13320
$ws.Range("M82").Value = -5294
$ws.Range("N82").Value = -14132
$ws.Range("H83").Value = 4333.3335
$ws.Range("I83").Value = 1150
$ws.Range("J83").Value = 4731.25
$ws.Range("K83").Value = 10350
$ws.Range("L83").Value = 42581.25
$ws.Range("M83").Value = -5670
$ws.Range("N83").Value = -51941.25
$ws.Range("H84").Value = 2197.6191
$ws.Range("J84").Value = 2271.0527
$ws.Range("L84").Value = 20439.4743
$ws.Range("N84").Value = -31671.4743
$ws.Range("H85").Value = 4016.6667
$ws.Range("I85").Value = 1900
$ws.Range("J85").Value = 4440
$ws.Range("K85").Value = 5700
$ws.Range("L85").Value = 13320
$ws.Range("M85").Value = -4296
$ws.Range("N85").Value = -16128
$ws.Range("H86").Value = 755.7143
$ws.Range("I86").Value = 398.57144
$ws.Range("J86").Value = 1112.8572
$ws.Range("K86").Value = 1195.71432
$ws.Range("L86").Value = 3338.5716
$ws.Range("M86").Value = -9.714320000000043
$ws.Range("N86").Value = -5710.571599999999
$ws.Range("H89").Value = 755.7143
$ws.Range("I89").Value = 398.57144
$ws.Range("J89").Value = 1112.8572
$ws.Range("K89").Value = 3587.14296
$ws.Range("L89").Value = 10015.7148
$ws.Range("M89").Value = 2340.85704
$ws.Range("N89").Value = -21871.7148
$ws.Range("H98").Value = 2761.7144
$ws.Range("I98").Value = 3333.3333
$ws.Range("J98").Value = 2333
$ws.Range("K98").Value = 9999.999899999999
$ws.Range("L98").Value = 6999
$ws.Range("M98").Value = -8501.999899999999
$ws.Range("N98").Value = -9995
$ws.Range("H122").Value = 613.3333
$ws.Range("I122").Value = 428.44446
$ws.Range("K122").Value = 3856.00014
$ws.Range("M122").Value = -1406.00014
$ws.Range("H131").Value = 16667693
$ws.Range("I131").Value = 71428990
$ws.Range("J131").Value = 1635181.1
$ws.Range("K131").Value = 214286970
$ws.Range("L131").Value = 4905543.300000001
$ws.Range("M131").Value = -214281930
$ws.Range("N131").Value = -4915623.300000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1942.8572
$ws.Range("I97").Value = 1942.8572
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1942.8572
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1446.8572
$ws.Range("N97").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5125.5
$ws.Range("I40").Value = 4667.3335
$ws.Range("K40").Value = 4667.3335
$ws.Range("M40").Value = -4531.3335
$ws.Range("H132").Value = 5115551.5
$ws.Range("I132").Value = 6668035.5
$ws.Range("J132").Value = 2528078.8
$ws.Range("K132").Value = 20004106.5
$ws.Range("L132").Value = 7584236.399999999
$ws.Range("M132").Value = -20001576.5
$ws.Range("N132").Value = -7589296.399999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 13650.75
$ws.Range("I96").Value = 17534.334
$ws.Range("J96").Value = 2000
$ws.Range("K96").Value = 17534.334
$ws.Range("L96").Value = 2000
$ws.Range("M96").Value = -16161.334
$ws.Range("N96").Value = -4746
$ws.Range("H122").Value = 21150.8
$ws.Range("I122").Value = 26063.5
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 78190.5
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -75740.5
$ws.Range("N122").Value = -9400
$ws.Range("H132").Value = 77909440
$ws.Range("I132").Value = 92309300
$ws.Range("J132").Value = 57109652
$ws.Range("K132").Value = 276927900
$ws.Range("L132").Value = 171328956
$ws.Range("M132").Value = -276925370
$ws.Range("N132").Value = -171334016

Write-Host "Applied all updates"
